# Insert two new daily price rows for "Perejil" (Primera / Segunda) at
# row 704, pushing the existing data block (old rows 704:789) down to
# 706:791. This grows the used range from A1:R789 to A1:R791.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("704:705").Insert()

# New row 704 - Perejil, Primera
$ws.Cells.Item(704, 1).Value  = 6
$ws.Cells.Item(704, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(704, 3).Value  = "Metropolitana"
$ws.Cells.Item(704, 4).Value  = 45142
$ws.Cells.Item(704, 5).Value  = 13
$ws.Cells.Item(704, 6).Value  = 100112044
$ws.Cells.Item(704, 7).Value  = "Perejil"
$ws.Cells.Item(704, 8).Value  = "Sin especificar"
$ws.Cells.Item(704, 9).Value  = "Primera"
$ws.Cells.Item(704, 10).Value = 150
$ws.Cells.Item(704, 11).Value = 13000
$ws.Cells.Item(704, 12).Value = 13000
$ws.Cells.Item(704, 13).Value = 13000
$ws.Cells.Item(704, 14).Value = "`$/docena de atados"
$ws.Cells.Item(704, 15).Value = "Región Metropolitana"
$ws.Cells.Item(704, 16).Value = 4333
$ws.Cells.Item(704, 17).Value = 3
$ws.Cells.Item(704, 18).Value = "Hortaliza"

# New row 705 - Perejil, Segunda
$ws.Cells.Item(705, 1).Value  = 6
$ws.Cells.Item(705, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(705, 3).Value  = "Metropolitana"
$ws.Cells.Item(705, 4).Value  = 45142
$ws.Cells.Item(705, 5).Value  = 13
$ws.Cells.Item(705, 6).Value  = 100112044
$ws.Cells.Item(705, 7).Value  = "Perejil"
$ws.Cells.Item(705, 8).Value  = "Sin especificar"
$ws.Cells.Item(705, 9).Value  = "Segunda"
$ws.Cells.Item(705, 10).Value = 180
$ws.Cells.Item(705, 11).Value = 12000
$ws.Cells.Item(705, 12).Value = 12000
$ws.Cells.Item(705, 13).Value = 12000
$ws.Cells.Item(705, 14).Value = "`$/docena de atados"
$ws.Cells.Item(705, 15).Value = "Región Metropolitana"
$ws.Cells.Item(705, 16).Value = 4000
$ws.Cells.Item(705, 17).Value = 3
$ws.Cells.Item(705, 18).Value = "Hortaliza"
